{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"210\u00f75=\", \"791\u00f76=\"],\n  [\"487\u00f74=\", \"757\u00f77=\"],\n  [\"709\u00f76=\", \"220\u00f78=\"],\n  [\"372\u00f78=\", \"350\u00f75=\"],\n  [\"430\u00f76=\", \"742\u00f76=\"],\n  [\"763\u00f74=\", \"770\u00f78=\"],\n  [\"903\u00f76=\", \"618\u00f77=\"],\n  [\"181\u00f76=\", \"524\u00f72=\"],\n  [\"138\u00f76=\", \"852\u00f74=\"],\n  [\"677\u00f72=\", \"711\u00f74=\"],\n  [\"155\u00f77=\", \"687\u00f73=\"],\n  [\"879\u00f79=\", \"845\u00f79=\"],\n  [\"353\u00f75=\", \"658\u00f77=\"],\n  [\"136\u00f75=\", \"468\u00f75=\"],\n  [\"493\u00f72=\", \"617\u00f75=\"],\n  [\"487\u00f78=\", \"534\u00f72=\"],\n  [\"341\u00f75=\", \"211\u00f72=\"],\n  [\"786\u00f73=\", \"605\u00f73=\"],\n  [\"164\u00f75=\", \"596\u00f75=\"],\n  [\"312\u00f77=\", \"639\u00f79=\"],\n  [\"754\u00f78=\", \"718\u00f75=\"],\n  [\"477\u00f72=\", \"444\u00f76=\"],\n  [\"159\u00f73=\", \"419\u00f75=\"],\n  [\"396\u00f77=\", \"112\u00f79=\"],\n  [\"179\u00f76=\", \"357\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"210\u00f75=\"; New=\"791\u00f76=\"}\n    @{Old=\"487\u00f74=\"; New=\"757\u00f77=\"}\n    @{Old=\"709\u00f76=\"; New=\"220\u00f78=\"}\n    @{Old=\"372\u00f78=\"; New=\"350\u00f75=\"}\n    @{Old=\"430\u00f76=\"; New=\"742\u00f76=\"}\n    @{Old=\"763\u00f74=\"; New=\"770\u00f78=\"}\n    @{Old=\"903\u00f76=\"; New=\"618\u00f77=\"}\n    @{Old=\"181\u00f76=\"; New=\"524\u00f72=\"}\n    @{Old=\"138\u00f76=\"; New=\"852\u00f74=\"}\n    @{Old=\"677\u00f72=\"; New=\"711\u00f74=\"}\n    @{Old=\"155\u00f77=\"; New=\"687\u00f73=\"}\n    @{Old=\"879\u00f79=\"; New=\"845\u00f79=\"}\n    @{Old=\"353\u00f75=\"; New=\"658\u00f77=\"}\n    @{Old=\"136\u00f75=\"; New=\"468\u00f75=\"}\n    @{Old=\"493\u00f72=\"; New=\"617\u00f75=\"}\n    @{Old=\"487\u00f78=\"; New=\"534\u00f72=\"}\n    @{Old=\"341\u00f75=\"; New=\"211\u00f72=\"}\n    @{Old=\"786\u00f73=\"; New=\"605\u00f73=\"}\n    @{Old=\"164\u00f75=\"; New=\"596\u00f75=\"}\n    @{Old=\"312\u00f77=\"; New=\"639\u00f79=\"}\n    @{Old=\"754\u00f78=\"; New=\"718\u00f75=\"}\n    @{Old=\"477\u00f72=\"; New=\"444\u00f76=\"}\n    @{Old=\"159\u00f73=\"; New=\"419\u00f75=\"}\n    @{Old=\"396\u00f77=\"; New=\"112\u00f79=\"}\n    @{Old=\"179\u00f76=\"; New=\"357\u00f74=\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
